# Update the "Database_DataDictionary" worksheet:
#  - Row 8 ("age_months") is renamed to "age_years" and its explanation is
#    updated to refer to years instead of months.
#  - Two new rows are appended describing the new "playcount" and
#    "time_since_firstgameplay" variables.
#  - The active selection is moved to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Variable names (column A) first ---
$ws.Range("A8").Value  = "age_years"
$ws.Range("A9").Value  = "playcount"
$ws.Range("A10").Value = "time_since_firstgameplay"

# --- Explanations (column B) second ---
$ws.Range("B8").Value  = "Calculated age in years based on the birthdate provided or approximated"
$ws.Range("B10").Value = "Time in months that have passed since a participant's first ACE session"
$ws.Range("B9").Value  = "Count of times the participant has played ACE, which may not correspond to the Time.Point in cases where participant missed a session (1-4) "

# Move the selection like the saved workbook shows.
$ws.Range("B13").Select()
